$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H19").Value = 1476.7059
$ws.Range("I19").Value = 1477.5834
$ws.Range("J19").Value = 1474.6
$ws.Range("K19").Value = 1477.5834
$ws.Range("L19").Value = 1474.6
$ws.Range("M19").Value = -1302.5834
$ws.Range("N19").Value = -1824.6
$ws.Range("H41").Value = 79.333336
$ws.Range("I41").Value = 70
$ws.Range("J41").Value = 84
$ws.Range("K41").Value = 70
$ws.Range("L41").Value = 84
$ws.Range("M41").Value = 370
$ws.Range("N41").Value = -964
$ws.Range("H100").Value = 3545.4473
$ws.Range("I100").Value = 2802.087
$ws.Range("K100").Value = 2802.087
$ws.Range("M100").Value = -2261.087
$ws.Range("H113").Value = 8304.417
$ws.Range("I113").Value = 8093.857
$ws.Range("J113").Value = 8599.2
$ws.Range("K113").Value = 8093.857
$ws.Range("L113").Value = 8599.2
$ws.Range("M113").Value = -4839.857
$ws.Range("N113").Value = -15107.2
$ws.Range("H116").Value = 20132.555
$ws.Range("J116").Value = 4439.6
$ws.Range("L116").Value = 4439.6
$ws.Range("N116").Value = -11323.6
$ws.Range("H129").Value = 1533.4375
$ws.Range("I129").Value = 1015.1111
$ws.Range("K129").Value = 3045.3333
$ws.Range("M129").Value = 1954.6667
$ws.Range("H141").Value = 5215.884
$ws.Range("I141").Value = 4233.657
$ws.Range("J141").Value = 9513.125
$ws.Range("K141").Value = 12700.971
$ws.Range("L141").Value = 28539.375
$ws.Range("M141").Value = -7520.971000000001
$ws.Range("N141").Value = -38899.375

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 17648826
$ws.Range("I32").Value = 17335486
$ws.Range("K32").Value = 17335486
$ws.Range("M32").Value = -17335199
$ws.Range("H61").Value = 3565.88
$ws.Range("I61").Value = 2842.1538
$ws.Range("K61").Value = 2842.1538
$ws.Range("M61").Value = -2630.1538
$ws.Range("H132").Value = 3080.8484
$ws.Range("I132").Value = 3061.8518
$ws.Range("J132").Value = 3166.3333
$ws.Range("K132").Value = 9185.5554
$ws.Range("L132").Value = 9498.999899999999
$ws.Range("M132").Value = -6655.555399999999
$ws.Range("N132").Value = -14558.9999
$ws.Range("H136").Value = 3565.88
$ws.Range("I136").Value = 2842.1538
$ws.Range("K136").Value = 8526.4614
$ws.Range("M136").Value = -5976.4614

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H7").Value = 14199.6
$ws.Range("J7").Value = 17000
$ws.Range("L7").Value = 17000
$ws.Range("N7").Value = -17226
$ws.Range("H134").Value = 9526622
$ws.Range("I134").Value = 9526622
$ws.Range("K134").Value = 28579866
$ws.Range("M134").Value = -28577331

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H7").Value = 24.583334
$ws.Range("I7").Value = 18.73913
$ws.Range("K7").Value = 18.73913
$ws.Range("M7").Value = 94.26087
$ws.Range("H31").Value = 4984.5776
$ws.Range("I31").Value = 2026.05
$ws.Range("K31").Value = 2026.05
$ws.Range("M31").Value = -1731.05
$ws.Range("H34").Value = 4984.5776
$ws.Range("I34").Value = 2026.05
$ws.Range("K34").Value = 2026.05
$ws.Range("M34").Value = -1824.05
$ws.Range("H58").Value = 2824.3274
$ws.Range("I58").Value = 2624.1277
$ws.Range("K58").Value = 2624.1277
$ws.Range("M58").Value = -2421.1277
$ws.Range("H60").Value = 52500
$ws.Range("J60").Value = 52500
$ws.Range("L60").Value = 52500
$ws.Range("N60").Value = -53522
$ws.Range("H68").Value = 0
$ws.Range("J68").Value = 0
$ws.Range("L68").Value = 0
$ws.Range("N68").ClearContents()
$ws.Range("H71").Value = 0
$ws.Range("J71").Value = 0
$ws.Range("L71").Value = 0
$ws.Range("N71").ClearContents()
$ws.Range("H98").Value = 61239.668
$ws.Range("J98").Value = 61239.668
$ws.Range("L98").Value = 61239.668
$ws.Range("N98").Value = -65731.668
$ws.Range("H99").Value = 2902.75
$ws.Range("I99").Value = 2937
$ws.Range("J99").Value = 2800
$ws.Range("K99").Value = 2937
$ws.Range("L99").Value = 2800
$ws.Range("M99").Value = -1439
$ws.Range("N99").Value = -5796
$ws.Range("H123").Value = 47124.75
$ws.Range("J123").Value = 47124.75
$ws.Range("L123").Value = 47124.75
$ws.Range("N123").Value = -56924.75
$ws.Range("H124").Value = 84999.5
$ws.Range("J124").Value = 84999.5
$ws.Range("L124").Value = 84999.5
$ws.Range("N124").Value = -89909.5
$ws.Range("H125").Value = 62331.332
$ws.Range("J125").Value = 62331.332
$ws.Range("L125").Value = 62331.332
$ws.Range("N125").Value = -67251.332
$ws.Range("H126").Value = 2902.75
$ws.Range("I126").Value = 2937
$ws.Range("J126").Value = 2800
$ws.Range("K126").Value = 8811
$ws.Range("L126").Value = 8400
$ws.Range("M126").Value = -6341
$ws.Range("N126").Value = -13340
$ws.Range("H130").Value = 47249.875
$ws.Range("J130").Value = 47249.875
$ws.Range("L130").Value = 47249.875
$ws.Range("N130").Value = -57289.875
$ws.Range("H131").Value = 5500
$ws.Range("J131").Value = 5500
$ws.Range("L131").Value = 5500
$ws.Range("N131").Value = -15580
$ws.Range("H132").Value = 3781.9756
$ws.Range("I132").Value = 3125.9443
$ws.Range("J132").Value = 8505.4
$ws.Range("K132").Value = 9377.832900000001
$ws.Range("L132").Value = 25516.2
$ws.Range("M132").Value = -6847.832900000001
$ws.Range("N132").Value = -30576.2
$ws.Range("H133").Value = 61865.332
$ws.Range("J133").Value = 61865.332
$ws.Range("L133").Value = 61865.332
$ws.Range("N133").Value = -66925.332
$ws.Range("H134").Value = 2746.8333
$ws.Range("I134").Value = 2158.125
$ws.Range("J134").Value = 3924.25
$ws.Range("K134").Value = 6474.375
$ws.Range("L134").Value = 11772.75
$ws.Range("M134").Value = -3939.375
$ws.Range("N134").Value = -16842.75
$ws.Range("H136").Value = 2824.3274
$ws.Range("I136").Value = 2624.1277
$ws.Range("K136").Value = 7872.3831
$ws.Range("M136").Value = -5322.3831

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H3").Value = 9684.667
$ws.Range("I3").Value = 9684.667
$ws.Range("K3").Value = 29054.001
$ws.Range("M3").Value = -28942.001
$ws.Range("H9").Value = 2502375
$ws.Range("I9").Value = 6667666.5
$ws.Range("J9").Value = 3200
$ws.Range("K9").Value = 20002999.5
$ws.Range("L9").Value = 9600
$ws.Range("M9").Value = -20002775.5
$ws.Range("N9").Value = -10048
$ws.Range("H23").Value = 190.91667
$ws.Range("J23").Value = 198.88889
$ws.Range("L23").Value = 596.6666700000001
$ws.Range("N23").Value = -1066.66667
$ws.Range("H39").Value = 4191.4546
$ws.Range("J39").Value = 4253.875
$ws.Range("L39").Value = 12761.625
$ws.Range("N39").Value = -13349.625
$ws.Range("H55").Value = 3978
$ws.Range("J55").Value = 5018.5713
$ws.Range("L55").Value = 15055.7139
$ws.Range("N55").Value = -15409.7139
$ws.Range("H133").Value = 4255.2
$ws.Range("I133").Value = 4041
$ws.Range("J133").Value = 4500
$ws.Range("K133").Value = 12123
$ws.Range("L133").Value = 13500
$ws.Range("M133").Value = -7063
$ws.Range("N133").Value = -23620
$ws.Range("H138").Value = 42601096
$ws.Range("I138").Value = 1823.3334
$ws.Range("K138").Value = 5470.0002
$ws.Range("M138").Value = -330.0002000000004
$ws.Range("H139").Value = 7999
$ws.Range("I139").Value = 0
$ws.Range("K139").Value = 0
$ws.Range("M139").ClearContents()

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 14158.333
$ws.Range("I70").Value = 77375
$ws.Range("K70").Value = 77375
$ws.Range("M70").Value = -77105
$ws.Range("H73").Value = 14158.333
$ws.Range("I73").Value = 77375
$ws.Range("K73").Value = 77375
$ws.Range("M73").Value = -76439
$ws.Range("H123").Value = 29592.334
$ws.Range("J123").Value = 29592.334
$ws.Range("L123").Value = 29592.334
$ws.Range("N123").Value = -34492.334
$ws.Range("H132").Value = 3399.55
$ws.Range("I132").Value = 3316.2942
$ws.Range("J132").Value = 3871.3333
$ws.Range("K132").Value = 9948.882599999999
$ws.Range("L132").Value = 11613.9999
$ws.Range("M132").Value = -7418.882599999999
$ws.Range("N132").Value = -16673.9999

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 7276.533
$ws.Range("I7").Value = 7158.3076
$ws.Range("K7").Value = 7158.3076
$ws.Range("M7").Value = -7046.3076
$ws.Range("H29").Value = 0
$ws.Range("I29").Value = 0
$ws.Range("K29").Value = 0
$ws.Range("M29").ClearContents()
$ws.Range("H126").Value = 7276.533
$ws.Range("I126").Value = 7158.3076
$ws.Range("K126").Value = 21474.9228
$ws.Range("M126").Value = -19004.9228

Write-Output "All changes applied"